# Update the "取得日時" (acquired datetime) column on the "ランサーズ" sheet
# from 2025-09-21 18:22:52 to 2025-09-21 18:30:27 for rows 2 through 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-09-21 18:22:52"
$newValue = "2025-09-21 18:30:27"

for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
